$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 74.2
$ws.Range("J9").Value = 25.25
$ws.Range("L9").Value = 25.25
$ws.Range("N9").Value = -363.25
$ws.Range("H26").Value = 3499.75
$ws.Range("J26").Value = 3499.75
$ws.Range("L26").Value = 3499.75
$ws.Range("N26").Value = -4187.75
$ws.Range("H113").Value = 5748.75
$ws.Range("I113").Value = 6500
$ws.Range("J113").Value = 4997.5
$ws.Range("K113").Value = 6500
$ws.Range("L113").Value = 4997.5
$ws.Range("M113").Value = -3246
$ws.Range("N113").Value = -11505.5
$ws.Range("H132").Value = 1016.0294
$ws.Range("J132").Value = 5249
$ws.Range("L132").Value = 15747
$ws.Range("N132").Value = -20807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6966.6333
$ws.Range("I32").Value = 4192.5
$ws.Range("K32").Value = 4192.5
$ws.Range("M32").Value = -3905.5
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H102").Value = 2473.8
$ws.Range("I102").Value = 1185
$ws.Range("K102").Value = 1185
$ws.Range("M102").Value = 437
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3443.6
$ws.Range("I134").Value = 2696
$ws.Range("K134").Value = 8088
$ws.Range("M134").Value = -5553

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4378.4443
$ws.Range("I31").Value = 3943
$ws.Range("K31").Value = 3943
$ws.Range("M31").Value = -3648
$ws.Range("H34").Value = 4378.4443
$ws.Range("I34").Value = 3943
$ws.Range("K34").Value = 3943
$ws.Range("M34").Value = -3741
$ws.Range("H86").Value = 6494
$ws.Range("I86").Value = 4237.5
$ws.Range("J86").Value = 9502.666999999999
$ws.Range("K86").Value = 4237.5
$ws.Range("L86").Value = 9502.666999999999
$ws.Range("M86").Value = -3114.5
$ws.Range("N86").Value = -11748.667
$ws.Range("H89").Value = 6494
$ws.Range("I89").Value = 4237.5
$ws.Range("J89").Value = 9502.666999999999
$ws.Range("K89").Value = 21187.5
$ws.Range("L89").Value = 47513.335
$ws.Range("M89").Value = -15571.5
$ws.Range("N89").Value = -58745.335
$ws.Range("H99").Value = 14569.857
$ws.Range("I99").Value = 11748
$ws.Range("J99").Value = 16306.385
$ws.Range("K99").Value = 11748
$ws.Range("L99").Value = 16306.385
$ws.Range("M99").Value = -10250
$ws.Range("N99").Value = -19302.385
$ws.Range("H126").Value = 14569.857
$ws.Range("I126").Value = 11748
$ws.Range("J126").Value = 16306.385
$ws.Range("K126").Value = 35244
$ws.Range("L126").Value = 48919.155
$ws.Range("M126").Value = -32774
$ws.Range("N126").Value = -53859.155
$ws.Range("H132").Value = 1459.2174
$ws.Range("I132").Value = 1492.4
$ws.Range("K132").Value = 4477.200000000001
$ws.Range("M132").Value = -1947.200000000001
$ws.Range("H134").Value = 2225.2122
$ws.Range("I134").Value = 2036.32
$ws.Range("K134").Value = 6108.96
$ws.Range("M134").Value = -3573.96

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 633.75
$ws.Range("J15").Value = 771.5
$ws.Range("L15").Value = 2314.5
$ws.Range("N15").Value = -2594.5
$ws.Range("H131").Value = 1333.2632
$ws.Range("I131").Value = 691.7143
$ws.Range("J131").Value = 1707.5
$ws.Range("K131").Value = 2075.1429
$ws.Range("L131").Value = 5122.5
$ws.Range("M131").Value = 2964.8571
$ws.Range("N131").Value = -15202.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 55801.39
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("N3").Value = -2232
$ws.Range("H35").Value = 778461.4399999999
$ws.Range("I35").Value = 778461.4399999999
$ws.Range("K35").Value = 778461.4399999999
$ws.Range("M35").Value = -778163.4399999999
$ws.Range("H122").Value = 46909.74
$ws.Range("I122").Value = 2947.25
$ws.Range("K122").Value = 8841.75
$ws.Range("M122").Value = -6391.75
$ws.Range("H132").Value = 2676.8928
$ws.Range("I132").Value = 2418.1052
$ws.Range("K132").Value = 7254.3156
$ws.Range("M132").Value = -4724.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5225
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224
$ws.Range("H39").Value = 59571.8
$ws.Range("I39").Value = 37859
$ws.Range("J39").Value = 65000
$ws.Range("K39").Value = 37859
$ws.Range("L39").Value = 65000
$ws.Range("M39").Value = -37399
$ws.Range("N39").Value = -65920
$ws.Range("H40").Value = 2956.5
$ws.Range("I40").Value = 2882.9167
$ws.Range("K40").Value = 2882.9167
$ws.Range("M40").Value = -2746.9167
$ws.Range("H126").Value = 5225
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 38500
$ws.Range("J48").Value = 38500
$ws.Range("L48").Value = 38500
$ws.Range("N48").Value = -39638
$ws.Range("H81").Value = 13692.714
$ws.Range("I81").Value = 10424.5
$ws.Range("J81").Value = 15000
$ws.Range("K81").Value = 20849
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = -19788
$ws.Range("N81").Value = -32122
$ws.Range("H84").Value = 13692.714
$ws.Range("I84").Value = 10424.5
$ws.Range("J84").Value = 15000
$ws.Range("K84").Value = 104245
$ws.Range("L84").Value = 150000
$ws.Range("M84").Value = -98941
$ws.Range("N84").Value = -160608
$ws.Range("H126").Value = 5666.5557
$ws.Range("I126").Value = 4833.6665
$ws.Range("K126").Value = 14500.9995
$ws.Range("M126").Value = -12030.9995
$ws.Range("H130").Value = 67498.5
$ws.Range("J130").Value = 67498.5
$ws.Range("L130").Value = 67498.5
$ws.Range("N130").Value = -77538.5
$ws.Range("H132").Value = 8767.223
$ws.Range("I132").Value = 1401
$ws.Range("J132").Value = 12450.333
$ws.Range("K132").Value = 4203
$ws.Range("L132").Value = 37350.999
$ws.Range("M132").Value = -1673
$ws.Range("N132").Value = -42410.999
